$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '60.913.41'
Set-TextValue $ws.Range("E2") '  +0.21%  '
Set-TextValue $ws.Range("D3") '2.352.24'
Set-TextValue $ws.Range("E3") '  -1.11%  '
Set-TextValue $ws.Range("E4") '  -0.01%  '
Set-TextValue $ws.Range("D5") '544.81'
Set-TextValue $ws.Range("E5") '  +0.03%  '
Set-TextValue $ws.Range("D6") '136.99'
Set-TextValue $ws.Range("E6") '  -3.00%  '
Set-TextValue $ws.Range("E7") '  -0.02%  '
Set-TextValue $ws.Range("D8") '0.525'
Set-TextValue $ws.Range("E8") '  -2.92%  '
Set-TextValue $ws.Range("D9") '2.350.78'
Set-TextValue $ws.Range("E9") '  -1.06%  '
Set-TextValue $ws.Range("E10") '  -0.18%  '
Set-TextValue $ws.Range("E11") '  +1.94%  '
Set-TextValue $ws.Range("E12") '  -0.44%  '
Set-TextValue $ws.Range("D13") '0.345'
Set-TextValue $ws.Range("E13") '  +0.08%  '
Set-TextValue $ws.Range("D14") '24.70'
Set-TextValue $ws.Range("E14") '  -3.16%  '
Set-TextValue $ws.Range("D15") '2.774.14'
Set-TextValue $ws.Range("E15") '  -1.13%  '
Set-TextValue $ws.Range("D16") '60.740.28'
Set-TextValue $ws.Range("E16") '  +0.18%  '
Set-TextValue $ws.Range("E17") '  -1.07%  '
Set-TextValue $ws.Range("D18") '2.351.59'
Set-TextValue $ws.Range("E18") '  -1.07%  '
Set-TextValue $ws.Range("D19") '10.64'
Set-TextValue $ws.Range("E19") '  +0.15%  '
Set-TextValue $ws.Range("D20") '320.18'
Set-TextValue $ws.Range("E20") '  +1.38%  '
Set-TextValue $ws.Range("D21") '4.13'
Set-TextValue $ws.Range("E21") '  +0.77%  '
Set-TextValue $ws.Range("D22") '6.58'
Set-TextValue $ws.Range("E22") '  -1.89%  '
Set-TextValue $ws.Range("E23") '  +0.03%  '
Set-TextValue $ws.Range("D24") '63.43'
Set-TextValue $ws.Range("E24") '  +1.07%  '
Set-TextValue $ws.Range("D25") '1.68'
Set-TextValue $ws.Range("E25") '  -7.36%  '
Set-TextValue $ws.Range("D26") '8.48'
Set-TextValue $ws.Range("E26") '  +9.62%  '
Set-TextValue $ws.Range("E27") '  +0.03%  '
Set-TextValue $ws.Range("B28") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D28") '8.00'
Set-TextValue $ws.Range("E28") '  +0.11%  '
Set-TextValue $ws.Range("B29") 'Fetch.AI'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D29") '1.38'
Set-TextValue $ws.Range("E29") '  -3.47%  '
Set-TextValue $ws.Range("B30") 'Bittensor'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D30") '498.01'
Set-TextValue $ws.Range("E30") '  -4.44%  '
Set-TextValue $ws.Range("B31") 'PEPE'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D31") '0.0₃0866'
Set-TextValue $ws.Range("E31") '  -6.70%  '
Set-TextValue $ws.Range("B32") 'Kaspa'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D32") '0.146'
Set-TextValue $ws.Range("E32") '  +1.73%  '
Set-TextValue $ws.Range("B33") 'PancakeSwap'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D33") '1.79'
Set-TextValue $ws.Range("E33") '  -2.37%  '
Set-TextValue $ws.Range("B34") 'ImmutableX'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D34") '1.51'
Set-TextValue $ws.Range("E34") '  -2.95%  '
Set-TextValue $ws.Range("B35") 'FirstDigitalUSD'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D35") '0.999'
Set-TextValue $ws.Range("E35") '  +0.03%  '
Set-TextValue $ws.Range("B36") 'NEARProtocol'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D36") '4.64'
Set-TextValue $ws.Range("E36") '  -0.29%  '
Set-TextValue $ws.Range("B37") 'PolygonEcosystemToken'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D37") '0.377'
Set-TextValue $ws.Range("E37") '  +0.39%  '
Set-TextValue $ws.Range("B38") 'EthereumClassic'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D38") '18.53'
Set-TextValue $ws.Range("E38") '  +2.91%  '
Set-TextValue $ws.Range("B39") 'Stacks'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D39") '1.83'
Set-TextValue $ws.Range("E39") '  +5.89%  '
Set-TextValue $ws.Range("B40") 'RenderToken'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D40") '5.25'
Set-TextValue $ws.Range("E40") '  -4.23%  '
Set-TextValue $ws.Range("B41") 'Monero'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D41") '144.33'
Set-TextValue $ws.Range("E41") '  +4.78%  '
Set-TextValue $ws.Range("B42") 'USDe'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D42") '0.999'
Set-TextValue $ws.Range("E42") '  -0.11%  '
Set-TextValue $ws.Range("B43") 'Aave'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D43") '142.78'
Set-TextValue $ws.Range("E43") '  +2.37%  '
Set-TextValue $ws.Range("B44") 'Filecoin'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D44") '3.57'
Set-TextValue $ws.Range("E44") '  +0.59%  '
Set-TextValue $ws.Range("B45") 'dogwifhat'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D45") '2.04'
Set-TextValue $ws.Range("E45") '  -8.03%  '
Set-TextValue $ws.Range("B46") 'Hedera'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D46") '0.0516'
Set-TextValue $ws.Range("E46") '  -0.25%  '
Set-TextValue $ws.Range("B47") 'InjectiveProtocol'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D47") '19.16'
Set-TextValue $ws.Range("E47") '  -5.50%  '
Set-TextValue $ws.Range("B48") 'Mantle'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D48") '0.568'
Set-TextValue $ws.Range("E48") '  -1.07%  '
Set-TextValue $ws.Range("B49") 'Stellar'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D49") '0.0901'
Set-TextValue $ws.Range("E49") '  -1.24%  '
Set-TextValue $ws.Range("B50") 'VeChain'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D50") '0.0221'
Set-TextValue $ws.Range("E50") '  -1.43%  '
Set-TextValue $ws.Range("B51") 'WhiteBITCoin'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D51") '11.39'
Set-TextValue $ws.Range("E51") '  +0.14%  '
